# Commit: "Added files for Photo and Test Report Test Cases"
#
# 1. PhotoTest: "run" -> "skip" on row 7 (choosefavoritephoto), and a new
#    row 8 for the "negativescenerios" test case.
# 2. A brand new "TestReportTest" sheet (after PhotoTest) with three new
#    test cases: uploadtestreport, uploadtestreportwithotheroption and
#    updatetestreport.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. PhotoTest updates
# ---------------------------------------------------------------------
$photo = $wb.Worksheets.Item("PhotoTest")

# choosefavoritephoto used to be "skip"ped via the "run" param - mark it
# "skip" explicitly and append the new negative-scenario test case.
$photo.Range("B7").Value = "skip"

$photo.Range("A8").Value = "negativescenerios"
$photo.Range("B8").Value = "skip"
$photo.Range("C8").Value = "Photo"

# Match the formatting already used by the "Params" column on this sheet.
$photo.Range("C7").Copy()
$photo.Range("C8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. New TestReportTest sheet (added after PhotoTest, becomes active)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "TestReportTest"

$newSheet.Range("A1").Value = "Test Case Name"
$newSheet.Range("B1").Value = "Run Status"
$newSheet.Range("C1").Value = "Params"

$newSheet.Range("A2").Value = "uploadtestreport"
$newSheet.Range("B2").Value = "skip"
$newSheet.Range("C2").Value = "Automation_Test,Vibhor,model,Test.jpeg,Test Description"

$newSheet.Range("A3").Value = "uploadtestreportwithotheroption"
$newSheet.Range("B3").Value = "skip"
$newSheet.Range("C3").Value = "Automation_Test,Vibhor,model,Test.jpeg,Test Description"

$newSheet.Range("A4").Value = "updatetestreport"
$newSheet.Range("B4").Value = "run"
$newSheet.Range("C4").Value = "Automation_Test,Vibhor,model,Updated Description,Test Description"

# Column widths matching the rest of the workbook's test sheets.
$newSheet.Columns.Item(1).ColumnWidth = 29.171875
$newSheet.Columns.Item(3).ColumnWidth = 32.2890625

# Reuse the existing header / params styling from PhotoTest.
$photo.Range("A1:C1").Copy()
$newSheet.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$photo.Range("C2").Copy()
$newSheet.Range("C2:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newSheet.Rows.Item(1).RowHeight = 22

# ---------------------------------------------------------------------
# 3. View state: selections on both sheets + TestReportTest as active tab
# ---------------------------------------------------------------------
$photo.Range("K25").Select()

$newSheet.Activate()
$newSheet.Range("N24").Select()
